$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-assign the "nome" values for rows 2-5
$ws.Range("B2").Value = "Rodrigo"
$ws.Range("B3").Value = "Paulo"
$ws.Range("B4").Value = "Abraão"
$ws.Range("B5").Value = "João"

# Delete row 6 (previously holding codigo 5 / "Jjj") entirely
$ws.Rows("6").Delete()
